$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily mark column D as text so numeric-looking price strings
# (e.g. "1.001") are stored verbatim instead of being parsed as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.844.45"
$ws.Range("E2").Value = "  -0.10%  "
$ws.Range("D3").Value = "1.887.89"
$ws.Range("E3").Value = "  -0.32%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "0.7466"
$ws.Range("E5").Value = "  -4.56%  "
$ws.Range("D6").Value = "242.39"
$ws.Range("E6").Value = "  -0.51%  "
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "0.3117"
$ws.Range("E8").Value = "  -0.60%  "
$ws.Range("D9").Value = "25.28"
$ws.Range("E9").Value = "  -1.79%  "
$ws.Range("D10").Value = "0.07131"
$ws.Range("E10").Value = "  -1.87%  "
$ws.Range("D11").Value = "0.08494"
$ws.Range("E11").Value = "  +4.94%  "
$ws.Range("D12").Value = "0.7596"
$ws.Range("E12").Value = "  -1.76%  "
$ws.Range("D13").Value = "1.888.98"
$ws.Range("E13").Value = "  -0.24%  "
$ws.Range("D14").Value = "5.355"
$ws.Range("E14").Value = "  -2.07%  "
$ws.Range("D15").Value = "93.35"
$ws.Range("E15").Value = "  -0.73%  "
$ws.Range("D16").Value = "6.174"
$ws.Range("E16").Value = "  -0.57%  "
$ws.Range("D17").Value = "29.863.25"
$ws.Range("E17").Value = "  +0.06%  "
$ws.Range("E18").Value = "  -1.64%  "
$ws.Range("D19").Value = "243.42"
$ws.Range("E19").Value = "  -1.39%  "
$ws.Range("D20").Value = "0.000007796"
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("D21").Value = "2.154.45"
$ws.Range("E21").Value = "  +1.83%  "
$ws.Range("D22").Value = "0.9997"
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").Value = "8.000"
$ws.Range("E23").Value = "  -1.45%  "
$ws.Range("D24").Value = "1.000"
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("D25").Value = "0.1587"
$ws.Range("E25").Value = "  -0.45%  "
$ws.Range("D26").Value = "9.355"
$ws.Range("E26").Value = "  -0.93%  "
$ws.Range("D27").Value = "162.43"
$ws.Range("E27").Value = "  -1.29%  "
$ws.Range("D28").Value = "18.74"
$ws.Range("E28").Value = "  -0.11%  "
$ws.Range("D29").Value = "2.025"
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("D30").Value = "1.494"
$ws.Range("E30").Value = "  +3.58%  "
$ws.Range("D31").Value = "1.531"
$ws.Range("E31").Value = "  -0.73%  "
$ws.Range("D32").Value = "4.495"
$ws.Range("E32").Value = "  +0.38%  "
$ws.Range("D33").Value = "4.111"
$ws.Range("E33").Value = "  +1.04%  "
$ws.Range("D34").Value = "0.05408"
$ws.Range("E34").Value = "  -2.87%  "
$ws.Range("D35").Value = "1.238"
$ws.Range("E35").Value = "  -0.23%  "
$ws.Range("D36").Value = "0.7450"
$ws.Range("E36").Value = "  -1.10%  "
$ws.Range("D37").Value = "1.003"
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("D38").Value = "2.710"
$ws.Range("E38").Value = "  +1.02%  "
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("D40").Value = "2.771"
$ws.Range("E40").Value = "  -0.95%  "
$ws.Range("D41").Value = "0.4449"
$ws.Range("E41").Value = "  -0.32%  "
$ws.Range("D42").Value = "6.074"
$ws.Range("E42").Value = "  +1.78%  "
$ws.Range("D43").Value = "1.090.54"
$ws.Range("E43").Value = "  -4.24%  "
$ws.Range("D44").Value = "72.46"
$ws.Range("E44").Value = "  -2.25%  "
$ws.Range("D45").Value = "0.8542"
$ws.Range("E45").Value = "  +0.24%  "
$ws.Range("D46").Value = "1.001"
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("D47").Value = "102.31"
$ws.Range("E47").Value = "  -0.06%  "
$ws.Range("D48").Value = "7.684"
$ws.Range("E48").Value = "  +1.99%  "
$ws.Range("D49").Value = "1.861"
$ws.Range("E49").Value = "  -1.49%  "
$ws.Range("D50").Value = "3.048"
$ws.Range("E50").Value = "  -2.68%  "
$ws.Range("D51").Value = "2.049.29"
$ws.Range("E51").Value = "  +0.07%  "

# Restore the original (default) cell style now that the text values are set.
$ws.Range("D2:D51").Style = "Normal"

